$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A59").Value = "0OMKJG"
$ws.Range("B59").Value = "Film de fusor HP"
$ws.Range("C59").Value = "P4010 P4014 P4015 P4510 P4515 M455 M600 M601 M602 M603 M604 M605 M606 M630"
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = 220000
$ws.Range("F59").Value = 3
$ws.Range("G59").Value = 0
$ws.Range("H59").Formula = "=(E59-D59)*G59"
$ws.Range("I59").Formula = "=D59*F59"
$ws.Range("J59").Value = 0
